$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "41.140.68"
$ws.Cells.Item(2, 5).Value = "  -1.64%  "

$ws.Cells.Item(3, 4).Value = "2.178.10"
$ws.Cells.Item(3, 5).Value = "  -2.23%  "

$ws.Cells.Item(4, 5).Value = "  -0.11%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "248.73"
$ws.Cells.Item(5, 5).Value = "  -0.35%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.615"
$ws.Cells.Item(6, 5).Value = "  -2.62%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "66.07"
$ws.Cells.Item(7, 5).Value = "  -7.87%  "

$ws.Cells.Item(8, 5).Value = "  -0.02%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.570"
$ws.Cells.Item(9, 5).Value = "  -3.74%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "58.88"
$ws.Cells.Item(10, 5).Value = "  +0.93%  "

$ws.Cells.Item(11, 2).Value = "Avalanche"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "35.97"
$ws.Cells.Item(11, 5).Value = "  -13.22%  "

$ws.Cells.Item(12, 2).Value = "Dogecoin"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0927"
$ws.Cells.Item(12, 5).Value = "  -4.95%  "

$ws.Cells.Item(13, 5).Value = "  -2.02%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.86"
$ws.Cells.Item(14, 5).Value = "  -4.89%  "

$ws.Cells.Item(15, 4).Value = "2.500.57"
$ws.Cells.Item(15, 5).Value = "  -2.33%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "14.37"
$ws.Cells.Item(16, 5).Value = "  -4.38%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.852"
$ws.Cells.Item(17, 5).Value = "  -1.68%  "

$ws.Cells.Item(18, 4).Value = "2.184.26"
$ws.Cells.Item(18, 5).Value = "  -1.85%  "

$ws.Cells.Item(19, 4).Value = "41.089.42"
$ws.Cells.Item(19, 5).Value = "  -1.69%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0942"
$ws.Cells.Item(20, 5).Value = "  -3.20%  "

$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.09"
$ws.Cells.Item(21, 5).Value = "  -2.30%  "

$ws.Cells.Item(22, 2).Value = "Litecoin"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "71.65"
$ws.Cells.Item(22, 5).Value = "  -2.11%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "228.28"
$ws.Cells.Item(23, 5).Value = "  -3.34%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.04"
$ws.Cells.Item(24, 5).Value = "  -6.07%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.77"
$ws.Cells.Item(25, 5).Value = "  -6.39%  "

$ws.Cells.Item(26, 5).Value = "  +0.05%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.26"
$ws.Cells.Item(27, 5).Value = "  +4.70%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.41"
$ws.Cells.Item(28, 5).Value = "  -5.76%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "3.72"
$ws.Cells.Item(29, 5).Value = "  -4.60%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "167.64"
$ws.Cells.Item(30, 5).Value = "  -2.41%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.01"
$ws.Cells.Item(31, 5).Value = "  -4.75%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "20.19"
$ws.Cells.Item(32, 5).Value = "  -3.43%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.122"
$ws.Cells.Item(33, 5).Value = "  -0.83%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.73"
$ws.Cells.Item(34, 5).Value = "  +2.08%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0736"
$ws.Cells.Item(35, 5).Value = "  +0.91%  "

$ws.Cells.Item(36, 5).Value = "  -3.47%  "

$ws.Cells.Item(37, 5).Value = "  -4.29%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.97"
$ws.Cells.Item(38, 5).Value = "  +0.12%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "24.49"
$ws.Cells.Item(39, 5).Value = "  -5.97%  "

$ws.Cells.Item(40, 5).Value = "  +0.51%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.21"
$ws.Cells.Item(41, 5).Value = "  -4.06%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.35"
$ws.Cells.Item(42, 5).Value = "  +10.04%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.47"
$ws.Cells.Item(43, 5).Value = "  -8.67%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "60.86"
$ws.Cells.Item(44, 5).Value = "  -10.93%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "11.29"
$ws.Cells.Item(45, 5).Value = "  -5.72%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "8.53"
$ws.Cells.Item(46, 5).Value = "  -2.67%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.190"
$ws.Cells.Item(47, 5).Value = "  -8.95%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0995"
$ws.Cells.Item(48, 5).Value = "  -2.54%  "

$ws.Cells.Item(49, 5).Value = "  -0.22%  "

$ws.Cells.Item(50, 5).Value = "  -1.27%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.15"
$ws.Cells.Item(51, 5).Value = "  -3.90%  "
